$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176612734794617
$ws.Range("B1").Value = 2.364291906356812
$ws.Range("C1").Value = 3.469253063201904
$ws.Range("D1").Value = 1.771709084510803
$ws.Range("E1").Value = 1.211493134498596
